$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated roster data (player, position, team) replacing the previous 17 rows
# and adding one new row (Ayo Dosunmu), for a total of 18 data rows (A2:C19).
$data = @(
    @("Chris Paul",          "PG",    "San Antonio Spurs"),
    @("Jalen Suggs",         "PG,SG", "Orlando Magic"),
    @("Bogdan Bogdanovic",   "SG,SF", "Atlanta Hawks"),
    @("Jaylen Brown",        "SG,SF", "Boston Celtics"),
    @("Julian Champagnie",   "SF,PF", "San Antonio Spurs"),
    @("Deni Avdija",         "SF,PF", "Portland Trail Blazers"),
    @("Pascal Siakam",       "SF,PF", "Indiana Pacers"),
    @("Jerami Grant",        "SF,PF", "Portland Trail Blazers"),
    @("Clint Capela",        "C",     "Atlanta Hawks"),
    @("Rudy Gobert",         "C",     "Minnesota Timberwolves"),
    @("Ayo Dosunmu",         "SG,SF", "Chicago Bulls"),
    @("Nikola Jokic",        "C",     "Denver Nuggets"),
    @("Dejounte Murray",     "PG,SG", "New Orleans Pelicans"),
    @("Russell Westbrook",   "PG",    "Denver Nuggets"),
    @("Jalen Green",         "PG,SG", "Houston Rockets"),
    @("Paolo Banchero",      "SF,PF", "Orlando Magic"),
    @("Chet Holmgren",       "PF,C",  "Oklahoma City Thunder"),
    @("Jakob Poeltl",        "C",     "Toronto Raptors")
)

$rowCount = $data.Count
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
